$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so they stay as literal text
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.336.90"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.844.47"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").Value = "0.9979"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "240.17"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "0.6270"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "0.9984"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "0.07483"
$ws.Range("E8").Value = "  -1.91%  "
$ws.Range("D9").Value = "0.2902"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").Value = "24.35"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").Value = "0.07707"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "1.844.59"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("D13").Value = "5.001"
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").Value = "0.6789"
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").Value = "0.00001025"
$ws.Range("E15").Value = "  -3.20%  "
$ws.Range("D16").Value = "82.14"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").Value = "2.101.46"
$ws.Range("E17").Value = "  -3.91%  "
$ws.Range("D18").Value = "6.157"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "29.367.82"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "228.86"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").Value = "12.33"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "0.9983"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").Value = "0.9983"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "158.70"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").Value = "0.1377"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").Value = "8.401"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "17.55"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").Value = "0.06387"
$ws.Range("E29").Value = "  +14.16%  "
$ws.Range("D30").Value = "1.382"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").Value = "1.473"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("D32").Value = "4.092"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").Value = "4.059"
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").Value = "1.823"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").Value = "1.140"
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("D36").Value = "0.6986"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "2.575"
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("D38").Value = "1.258.30"
$ws.Range("E38").Value = "  +2.34%  "
$ws.Range("D39").Value = "2.826"
$ws.Range("E39").Value = "  +3.83%  "
$ws.Range("D40").Value = "0.01829"
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("D41").Value = "6.586"
$ws.Range("E41").Value = "  +3.12%  "
$ws.Range("D42").Value = "0.9078"
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").Value = "0.9976"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").Value = "2.007.38"
$ws.Range("E44").Value = "  -18.43%  "
$ws.Range("D45").Value = "101.50"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").Value = "66.07"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "0.00000000118"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.723"
$ws.Range("E48").Value = "  +2.61%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.1176"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "7.063"
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").Value = "9.047"
$ws.Range("E51").Value = "  +0.62%  "
